# This script updates the "cryptos" price/volume listing on Sheet1 to
# reflect the refreshed scrape values from the GitHub Actions run.
# Two coin pairs also swapped rank position (rows 15/16, 36/37, 50/51),
# so their Coin name / Link / Price / Volume cells are fully replaced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new literal value. Using a hashtable keeps the
# long list of per-cell updates easy to scan / maintain.
$updates = [ordered]@{
    'D2' = '68.888.69'
    'E2' = '  -3.40%  '
    'D3' = '3.491.38'
    'E3' = '  -5.34%  '
    'E4' = '  -0.06%  '
    'D5' = '577.89'
    'E5' = '  -1.16%  '
    'D6' = '170.89'
    'E6' = '  -4.81%  '
    'D7' = '0.608'
    'E7' = '  -0.90%  '
    'D8' = '3.485.27'
    'E8' = '  -5.31%  '
    'E9' = '  +0.03%  '
    'D10' = '0.188'
    'E10' = '  -6.35%  '
    'D11' = '6.71'
    'E11' = '  +3.58%  '
    'D12' = '0.585'
    'E12' = '  -4.27%  '
    'D13' = '46.39'
    'E13' = '  -6.26%  '
    'D14' = '0.0000272'
    'E14' = '  -4.86%  '
    'B15' = 'WrappedliquidstakedEther2.0'
    'C15' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'D15' = '4.043.49'
    'E15' = '  -5.57%  '
    'B16' = 'BitcoinCash'
    'C16' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'D16' = '644.48'
    'E16' = '  -5.74%  '
    'D17' = '8.54'
    'E17' = '  -5.37%  '
    'D18' = '68.791.37'
    'E18' = '  -3.74%  '
    'D19' = '3.487.71'
    'E19' = '  -5.14%  '
    'E20' = '  -1.21%  '
    'D21' = '17.29'
    'E21' = '  -4.01%  '
    'D22' = '11.06'
    'E22' = '  -4.86%  '
    'D23' = '0.885'
    'E23' = '  -5.75%  '
    'D24' = '15.97'
    'E24' = '  -8.68%  '
    'D25' = '97.05'
    'E25' = '  -5.04%  '
    'D26' = '3.81'
    'E26' = '  -4.87%  '
    'E27' = '  +0.00%  '
    'D28' = '2.63'
    'E28' = '  -7.59%  '
    'D29' = '9.31'
    'E29' = '  -9.99%  '
    'D30' = '32.56'
    'E30' = '  -7.78%  '
    'D31' = '3.17'
    'E31' = '  -8.57%  '
    'D32' = '8.49'
    'E32' = '  -7.63%  '
    'D33' = '1.32'
    'E33' = '  -8.86%  '
    'D34' = '7.07'
    'E34' = '  -3.90%  '
    'D35' = '617.07'
    'E35' = '  +5.74%  '
    'B36' = 'Cosmos'
    'C36' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D36' = '10.77'
    'E36' = '  -4.22%  '
    'B37' = 'dogwifhat'
    'C37' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D37' = '3.52'
    'E37' = '  -14.98%  '
    'D38' = '0.103'
    'E38' = '  -5.09%  '
    'D39' = '56.59'
    'E39' = '  -4.10%  '
    'D40' = '0.999'
    'E40' = '  +0.08%  '
    'D41' = '0.136'
    'E41' = '  -6.93%  '
    'D42' = '0.0434'
    'E42' = '  -5.87%  '
    'D43' = '0.329'
    'E43' = '  -5.68%  '
    'D44' = '3.346.51'
    'E44' = '  -8.96%  '
    'D45' = '32.90'
    'E45' = '  -7.70%  '
    'D46' = '0.0₃0698'
    'E46' = '  -9.66%  '
    'D47' = '2.57'
    'E47' = '  -8.51%  '
    'D48' = '2.77'
    'E48' = '  -4.62%  '
    'D49' = '0.131'
    'E49' = '  -2.42%  '
    'B50' = 'MXToken'
    'C50' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D50' = '5.78'
    'E50' = '  +16.96%  '
    'B51' = 'Monero'
    'C51' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D51' = '132.94'
    'E51' = '  -3.00%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force Text format *before* writing the value so Excel does not
    # reinterpret numeric-looking strings (e.g. "577.89", "0.999") or
    # percent-looking strings (e.g. "  -3.40%  ") as numbers/dates.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
